# Generate Report for Handoff
# This script updates the localization-status workbook to reflect a fresh
# handoff-xliff generation run for the "0e7852e7-b2b4-435e-9d80-896e17311fc2"
# item: its priority flips from "low" to "ht", and its handoff timestamps
# are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# 0e7852e7 file (rows 4-7 all share this value via the same underlying text).
$wsOverview.Range("G4:G7").Value = "2016-08-19 00:31:32"

# zh-cn sheet: Priority (E) low -> ht, and Latest Handoff Datetime (H) refreshed
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-08-19 00:31:27"
$wsZhCn.Range("H5").Value = "2016-08-19 00:31:27"
$wsZhCn.Range("H6").Value = "2016-08-19 00:31:27"
$wsZhCn.Range("H7").Value = "2016-08-19 00:31:27"

# de-de sheet: Priority (E) low -> ht
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

# de-de sheet: Latest Handoff Datetime (H) for 0e7852e7 file shares the same
# underlying text as Overview's G column, so it is refreshed together.
$wsDeDe.Range("H4:H7").Value = "2016-08-19 00:31:32"
